# Update header row (row 1) text on the active worksheet ("Tabelle1")
# to reflect renamed column headers (EI NOx units changed and the
# "EI inverse range" columns renamed to "flown distance").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A1").Value = "altitude [ft]"
$ws.Range("B1").Value = "pressure level [hPa]"
$ws.Range("C1").Value = "EI NOx regional [g/kg(Fuel)]"
$ws.Range("D1").Value = "EI NOx single-aisle [g/kg(Fuel)]"
$ws.Range("E1").Value = "EI NOx wide-body [g/kg(Fuel)]"
$ws.Range("F1").Value = "flown distance [km/kg(fuel)] regional"
$ws.Range("G1").Value = "flown distance [km/kg(fuel)] single-aisle"
$ws.Range("H1").Value = "flown distance [km/kg(fuel)] wide-body"

# Update the selected/active cell on the sheet to match the saved view.
$ws.Range("H10").Select()
